# Updated cryptos list on Fri Feb 24 21:57:46 UTC 2023 with GitHub Actions
#
# Refresh the Price (D) and Volume(1h) (E) columns for every coin row, and
# swap Algorand/InternetComputer(DFINITY) which changed rank (rows 37/38).
#
# Note: several Price values (e.g. "1.003", "0.3752", "0.00001260") look
# like plain decimal numbers to Excel's automatic type detection and would
# otherwise be silently coerced to numeric cells (losing significant
# trailing zeros). They are written with a leading apostrophe to force text
# interpretation, matching the original inline-string cell type, and the
# cell style is reset back to "Normal" afterwards so no stray number format
# is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.132.21"
$ws.Range("E2").Value = "  -3.41%  "
$ws.Range("D3").Value = "1.604.67"
$ws.Range("E3").Value = "  -2.75%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'1.002"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "'301.84"
$ws.Range("E6").Value = "  -2.45%  "
$ws.Range("D7").Value = "'0.3752"
$ws.Range("D8").Value = "'0.3622"
$ws.Range("E8").Value = "  -5.38%  "
$ws.Range("D9").Value = "'48.59"
$ws.Range("E9").Value = "  -5.10%  "
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").Value = "'1.257"
$ws.Range("E11").Value = "  -6.73%  "
$ws.Range("D12").Value = "'0.08036"
$ws.Range("E12").Value = "  -4.67%  "
$ws.Range("D13").Value = "'22.82"
$ws.Range("E13").Value = "  -4.39%  "
$ws.Range("D14").Value = "'6.528"
$ws.Range("D15").Value = "'7.640"
$ws.Range("E15").Value = "  -2.93%  "
$ws.Range("D16").Value = "'0.00001260"
$ws.Range("E16").Value = "  -4.17%  "
$ws.Range("D17").Value = "1.606.45"
$ws.Range("E17").Value = "  -2.84%  "
$ws.Range("D18").Value = "'91.24"
$ws.Range("E18").Value = "  -3.33%  "
$ws.Range("D19").Value = "'0.06771"
$ws.Range("E19").Value = "  -3.20%  "
$ws.Range("D20").Value = "'18.25"
$ws.Range("E20").Value = "  -7.38%  "
$ws.Range("D21").Value = "'6.540"
$ws.Range("E21").Value = "  -5.60%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'13.03"
$ws.Range("E23").Value = "  -4.95%  "
$ws.Range("D24").Value = "23.145.04"
$ws.Range("E24").Value = "  -3.37%  "
$ws.Range("D25").Value = "'2.341"
$ws.Range("E25").Value = "  -4.70%  "
$ws.Range("D26").Value = "'2.881"
$ws.Range("E26").Value = "  -3.16%  "
$ws.Range("D27").Value = "'20.98"
$ws.Range("E27").Value = "  -4.86%  "
$ws.Range("D28").Value = "'150.13"
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").Value = "'5.254"
$ws.Range("E29").Value = "  -2.98%  "
$ws.Range("D30").Value = "'132.07"
$ws.Range("E30").Value = "  -4.73%  "
$ws.Range("D31").Value = "'2.389"
$ws.Range("E31").Value = "  -4.13%  "
$ws.Range("D32").Value = "'6.686"
$ws.Range("E32").Value = "  -14.34%  "
$ws.Range("D33").Value = "1.783.45"
$ws.Range("E33").Value = "  -2.80%  "
$ws.Range("D34").Value = "'0.9662"
$ws.Range("E34").Value = "  -7.68%  "
$ws.Range("D35").Value = "'0.07690"
$ws.Range("E35").Value = "  -4.83%  "
$ws.Range("D36").Value = "'0.02746"
$ws.Range("E36").Value = "  -7.22%  "
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "'0.2526"
$ws.Range("E37").Value = "  -5.90%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'6.172"
$ws.Range("E38").Value = "  -8.24%  "
$ws.Range("D39").Value = "'0.08825"
$ws.Range("E39").Value = "  -3.30%  "
$ws.Range("D40").Value = "'10.03"
$ws.Range("E40").Value = "  -7.75%  "
$ws.Range("D41").Value = "'1.391"
$ws.Range("E41").Value = "  -2.14%  "
$ws.Range("D42").Value = "'0.7102"
$ws.Range("E42").Value = "  -5.92%  "
$ws.Range("D43").Value = "'12.70"
$ws.Range("D44").Value = "'15.76"
$ws.Range("E44").Value = "  -3.66%  "
$ws.Range("D45").Value = "'0.6556"
$ws.Range("E45").Value = "  -5.53%  "
$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "'2.278"
$ws.Range("E47").Value = "  -7.30%  "
$ws.Range("D48").Value = "'3.970"
$ws.Range("E48").Value = "  -2.88%  "
$ws.Range("D49").Value = "'0.07981"
$ws.Range("E49").Value = "  -3.52%  "
$ws.Range("D50").Value = "'131.01"
$ws.Range("E50").Value = "  -2.46%  "
$ws.Range("D51").Value = "'1.161"
$ws.Range("E51").Value = "  -3.52%  "

# Reset style on cells that needed a leading apostrophe to force text
# interpretation (Excel would otherwise read these as numbers), so the
# resulting cell keeps the default (unstyled) format like the rest of the sheet.
$resetCells = @("D4","D5","D6","D7","D8","D9","D11","D12","D13","D14","D15","D16","D18","D19","D20","D21","D22","D23","D25","D26","D27","D28","D29","D30","D31","D32","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $resetCells) {
    $ws.Range($addr).Style = "Normal"
}